# Merge the split <id>...</id> runs back into single runs for both
# occurrences in the document (p050r_1 and p050r_2), matching the
# "add newly downloaded tc, tcn, tl" commit that collapsed the
# per-fragment runs produced by the earlier TEI re-import.

$d = $word.ActiveDocument

$d.Content.Find.Execute("<id>p050r_1</id>", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p050r_1</id>", 2) | Out-Null

$d.Content.Find.Execute("<id>p050r_2</id>", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p050r_2</id>", 2) | Out-Null
